$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.876.57"
$ws.Range("E2").Value = "'  +0.85%  "

$ws.Range("D3").Value = "'1.816.56"
$ws.Range("E3").Value = "'  +1.54%  "

$ws.Range("D4").Value = "'1.002"
$ws.Range("E4").Value = "'  +0.08%  "

$ws.Range("D5").Value = "'309.28"
$ws.Range("E5").Value = "'  +0.46%  "

$ws.Range("E6").Value = "'  +0.03%  "

$ws.Range("D7").Value = "'0.4683"
$ws.Range("E7").Value = "'  +2.74%  "

$ws.Range("D8").Value = "'0.3686"
$ws.Range("E8").Value = "'  -0.20%  "

$ws.Range("D9").Value = "'0.07372"
$ws.Range("E9").Value = "'  +2.39%  "

$ws.Range("D10").Value = "'0.8702"
$ws.Range("E10").Value = "'  +1.96%  "

$ws.Range("E11").Value = "'  +0.29%  "

$ws.Range("D12").Value = "'1.749.39"
$ws.Range("E12").Value = "'  -2.21%  "

$ws.Range("E13").Value = "'  +1.56%  "

$ws.Range("D14").Value = "'92.34"
$ws.Range("E14").Value = "'  +2.47%  "

$ws.Range("D15").Value = "'0.07064"
$ws.Range("E15").Value = "'  +0.51%  "

$ws.Range("D16").Value = "'6.491"
$ws.Range("E16").Value = "'  +0.53%  "

$ws.Range("D17").Value = "'1.003"
$ws.Range("E17").Value = "'  +0.04%  "

$ws.Range("D18").Value = "'0.000008705"
$ws.Range("E18").Value = "'  +1.28%  "

$ws.Range("E19").Value = "'  +0.02%  "

$ws.Range("D20").Value = "'14.75"
$ws.Range("E20").Value = "'  +1.50%  "

$ws.Range("D21").Value = "'26.952.61"
$ws.Range("E21").Value = "'  +1.11%  "

$ws.Range("D22").Value = "'5.347"
$ws.Range("E22").Value = "'  +1.50%  "

$ws.Range("D23").Value = "'10.56"
$ws.Range("E23").Value = "'  +0.04%  "

$ws.Range("D24").Value = "'2.083.56"
$ws.Range("E24").Value = "'  +3.72%  "

$ws.Range("D25").Value = "'1.903"
$ws.Range("E25").Value = "'  -0.16%  "

$ws.Range("D26").Value = "'151.34"
$ws.Range("E26").Value = "'  +1.18%  "

$ws.Range("B27").Value = "'EthereumClassic"
$ws.Range("C27").Value = "'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").Value = "'18.36"
$ws.Range("E27").Value = "'  +1.71%  "

$ws.Range("B28").Value = "'LidoDAOToken"
$ws.Range("C28").Value = "'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D28").Value = "'2.169"
$ws.Range("E28").Value = "'  +1.85%  "

$ws.Range("D29").Value = "'5.324"
$ws.Range("E29").Value = "'  +2.55%  "

$ws.Range("D30").Value = "'115.70"
$ws.Range("E30").Value = "'  +1.72%  "

$ws.Range("D31").Value = "'0.08931"
$ws.Range("E31").Value = "'  +1.32%  "

$ws.Range("D32").Value = "'0.7680"
$ws.Range("E32").Value = "'  +1.97%  "

$ws.Range("D33").Value = "'1.162"
$ws.Range("E33").Value = "'  +0.79%  "

$ws.Range("D34").Value = "'4.507"
$ws.Range("E34").Value = "'  +1.72%  "

$ws.Range("D35").Value = "'2.902"
$ws.Range("E35").Value = "'  +0.68%  "

$ws.Range("E37").Value = "'  -1.56%  "

$ws.Range("D38").Value = "'0.01965"
$ws.Range("E38").Value = "'  +1.68%  "

$ws.Range("D39").Value = "'0.05282"
$ws.Range("E39").Value = "'  +1.74%  "

$ws.Range("D40").Value = "'2.943"
$ws.Range("E40").Value = "'  +1.96%  "

$ws.Range("D41").Value = "'7.267"
$ws.Range("E41").Value = "'  +2.20%  "

$ws.Range("D42").Value = "'0.5319"
$ws.Range("E42").Value = "'  +2.42%  "

$ws.Range("D43").Value = "'2.353"
$ws.Range("E43").Value = "'  +1.31%  "

$ws.Range("D44").Value = "'0.1667"
$ws.Range("E44").Value = "'  +1.82%  "

$ws.Range("D45").Value = "'8.423"
$ws.Range("E45").Value = "'  -0.26%  "

$ws.Range("D46").Value = "'0.4940"
$ws.Range("E46").Value = "'  +0.07%  "

$ws.Range("E47").Value = "'  +2.41%  "

$ws.Range("E48").Value = "'  +0.08%  "

$ws.Range("D49").Value = "'103.96"
$ws.Range("E49").Value = "'  +0.20%  "

$ws.Range("E50").Value = "'  +1.88%  "

$ws.Range("E51").Value = "'  +0.19%  "
